$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.842.18"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.416.94"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.418.76"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.24%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.118"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("D13").Value = "4.003.81"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("E16").Value = "  -7.63%  "
$ws.Range("D17").Value = "63.918.59"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "3.399.64"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -5.67%  "
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.08%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -5.55%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("E36").Value = "  -5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.76"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  +10.49%  "
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").Value = "2.791.95"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0722"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -7.12%  "
$ws.Range("E46").Value = "  -5.27%  "
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "330.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("E51").Value = "  -5.24%  "
